$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("O").Insert()
$ws.Range("O1").Value = "DOCREF4"
$ws.Range("D5").Select()
